# Apply cryptocurrency price/volume updates to Sheet1 (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.028.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.172.76'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.91%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('E6').Value = '  -2.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '67.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -6.51%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.567'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0925'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.20%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '35.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -16.38%  '
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.498.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.857'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.184.96'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '40.974.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('E20').Value = '  -3.47%  '
$ws.Range('E21').Value = '  -2.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('E24').Value = '  -8.60%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.03%  '
$ws.Range('E28').Value = '  -3.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.91'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('E30').Value = '  -9.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.121'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.70'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.92%  '
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('E35').Value = '  -3.67%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '25.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.21%  '
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0300'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.20%  '
$ws.Range('E41').Value = '  -9.77%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '61.36'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -11.90%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.27%  '
$ws.Range('E45').Value = '  -11.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.22%  '
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0981'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.16%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.17'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('E51').Value = '  -0.41%  '
